# Saved progress at the end of the loop
# Update the "Qty executed upto date" (C) and the derived "Upto date Amount"
# (G / H) columns on the Bill Summary sheet with the latest measured
# quantities and their resulting amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Qty executed upto date (plain numbers) ---------------------------------
$ws.Range("C8").Value  = 55
$ws.Range("C9").Value  = 1
$ws.Range("C10").Value = 61
$ws.Range("C11").Value = 58
$ws.Range("C12").Value = 32
$ws.Range("C13").Value = 23
$ws.Range("C14").Value = 70
$ws.Range("C15").Value = 26
$ws.Range("C16").Value = 6
$ws.Range("C17").Value = 68

# --- Upto date Amount cells are stored as literal text, e.g. "256.00" ------
# A plain `.Value = "256.00"` assignment gets auto-coerced to the number 256
# by Excel, so instead write it as a literal-string formula and immediately
# collapse it to its value (copy / paste-special values-only) - the normal
# COM trick to land a genuinely text typed result without leaving a live
# formula or touching the cell's number format / style.
function Set-TextAmount($address, $text) {
    $cell = $ws.Range($address)
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

Set-TextAmount "G9"  "256.00"
Set-TextAmount "G10" "28792.00"
Set-TextAmount "G11" "38396.00"
Set-TextAmount "G13" "3128.00"
Set-TextAmount "G14" "1610.00"
Set-TextAmount "G19" "72182.00"
Set-TextAmount "H19" "72182.00"
Set-TextAmount "G21" "72182.00"
Set-TextAmount "H21" "72182.00"
